$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row at position 18, pushing existing rows 18-37 down to 19-38.
$ws.Rows.Item(18).Insert()

# The newest entry (row 17) gets a refreshed September timestamp.
$ws.Cells.Item(17, 19).Value = "2024-09-03 07:57:20"

# The newly-inserted row 18 captures the entry's previous state (same details,
# prior timestamp) that used to occupy row 17 before the refresh.
$ws.Cells.Item(18, 18).Value = "electricity avoid disconnection tangedco"
$ws.Cells.Item(18, 19).Value = "2024-09-03 07:56:18"
